$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE (J2) moves from "001" to "002" - keep as text (leading zero
# must survive), K2 (REPORT_TYPE_CODE) is untouched since its value ("001")
# does not change.
$ws.Range("J2").Value = "'002"
$ws.Range("J2").Style = "Normal"

# NOTICE_DATE / REPORT_DATE text updates
$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric metric updates
$ws.Range("O2").Value = 981044756.5599999
$ws.Range("P2").Value = 137813839.79
$ws.Range("Q2").Value = 121597320.37
$ws.Range("R2").Value = 104.3060627933
$ws.Range("S2").Value = 304289358.29
$ws.Range("T2").Value = -27.3850741114
$ws.Range("U2").Value = 82460862.69
$ws.Range("V2").Value = 5.9596473445
$ws.Range("W2").Value = 284700105.44
$ws.Range("X2").Value = 157234960.58
$ws.Range("Y2").Value = -12.1069287865

# ADVANCE_RECEIVABLES / ADVANCE_RECEIVABLES_RATIO no longer reported
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()

$ws.Range("AB2").Value = 696344651.12
$ws.Range("AC2").Value = 16.2189328895
$ws.Range("AD2").Value = -6.0375218332
$ws.Range("AE2").Value = -36.0102432711
$ws.Range("AF2").Value = 308.4641788622
$ws.Range("AG2").Value = 29.0200934806
